$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1805555555555556
$ws.Range("C2").Value = 0.5740740740740741
$ws.Range("J2").Value = 0.004629629629629629
$ws.Range("P2").Value = 0.1388888888888889
$ws.Range("S2").Value = 0.1018518518518518
$ws.Range("C3").Value = 0.01587301587301587
$ws.Range("J3").Value = 0.03174603174603174
$ws.Range("P3").Value = 0.7222222222222222
$ws.Range("S3").Value = 0.2301587301587301
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.7241379310344828
$ws.Range("S4").Value = 0.2068965517241379
$ws.Range("B6").Value = 0.06842105263157895
$ws.Range("F6").Value = 0.06842105263157895
$ws.Range("J6").Value = 0.2421052631578947
$ws.Range("O6").Value = 0.04210526315789474
$ws.Range("Q6").Value = 0.1210526315789474
$ws.Range("R6").Value = 0.05789473684210526
$ws.Range("S6").Value = 0.4
$ws.Range("B7").Value = 0.09677419354838709
$ws.Range("D7").Value = 0.01612903225806452
$ws.Range("J7").Value = 0.1505376344086022
$ws.Range("O7").Value = 0.01075268817204301
$ws.Range("Q7").Value = 0.1290322580645161
$ws.Range("R7").Value = 0.05913978494623656
$ws.Range("S7").Value = 0.489247311827957
$ws.Range("B8").Value = 0.07563025210084033
$ws.Range("D8").Value = 0.008403361344537815
$ws.Range("E8").Value = 0.002801120448179272
$ws.Range("F8").Value = 0.05602240896358544
$ws.Range("J8").Value = 0.1316526610644258
$ws.Range("O8").Value = 0.02521008403361345
$ws.Range("Q8").Value = 0.15406162464986
$ws.Range("R8").Value = 0.07002801120448179
$ws.Range("S8").Value = 0.4761904761904762
$ws.Range("B9").Value = 0.08125
$ws.Range("D9").Value = 0.0125
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.06875000000000001
$ws.Range("O9").Value = 0.0125
$ws.Range("Q9").Value = 0.18125
$ws.Range("R9").Value = 0.10625
$ws.Range("S9").Value = 0.4875
$ws.Range("B10").Value = 0.09744094488188976
$ws.Range("D10").Value = 0.01968503937007874
$ws.Range("F10").Value = 0.06988188976377953
$ws.Range("J10").Value = 0.1141732283464567
$ws.Range("O10").Value = 0.01377952755905512
$ws.Range("Q10").Value = 0.187992125984252
$ws.Range("R10").Value = 0.03838582677165354
$ws.Range("S10").Value = 0.4586614173228347
$ws.Range("G11").Value = 0.1459074733096085
$ws.Range("J11").Value = 0.06761565836298933
$ws.Range("K11").Value = 0.1886120996441281
$ws.Range("L11").Value = 0.5729537366548043
$ws.Range("S11").Value = 0.02491103202846975
$ws.Range("G12").Value = 0.7668711656441718
$ws.Range("J12").Value = 0.1717791411042945
$ws.Range("K12").Value = 0.006134969325153374
$ws.Range("L12").Value = 0.01840490797546012
$ws.Range("S12").Value = 0.03680981595092025
$ws.Range("G13").Value = 0.84375
$ws.Range("S13").Value = 0.03125
$ws.Range("F15").Value = 0.02824858757062147
$ws.Range("H15").Value = 0.1412429378531073
$ws.Range("I15").Value = 0.1073446327683616
$ws.Range("J15").Value = 0.3163841807909605
$ws.Range("K15").Value = 0.07909604519774012
$ws.Range("M15").Value = 0.01694915254237288
$ws.Range("O15").Value = 0.05649717514124294
$ws.Range("S15").Value = 0.2542372881355932
$ws.Range("F16").Value = 0.02127659574468085
$ws.Range("H16").Value = 0.1276595744680851
$ws.Range("I16").Value = 0.09929078014184398
$ws.Range("J16").Value = 0.4397163120567376
$ws.Range("K16").Value = 0.09929078014184398
$ws.Range("O16").Value = 0.0425531914893617
$ws.Range("S16").Value = 0.1702127659574468
$ws.Range("F17").Value = 0.03174603174603174
$ws.Range("H17").Value = 0.1396825396825397
$ws.Range("I17").Value = 0.09206349206349207
$ws.Range("J17").Value = 0.4222222222222222
$ws.Range("K17").Value = 0.09841269841269841
$ws.Range("M17").Value = 0.01587301587301587
$ws.Range("O17").Value = 0.06666666666666667
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("F18").Value = 0.009803921568627451
$ws.Range("H18").Value = 0.2254901960784314
$ws.Range("I18").Value = 0.07843137254901961
$ws.Range("J18").Value = 0.392156862745098
$ws.Range("K18").Value = 0.1274509803921569
$ws.Range("M18").Value = 0.0196078431372549
$ws.Range("O18").Value = 0.0392156862745098
$ws.Range("S18").Value = 0.107843137254902
$ws.Range("F19").Value = 0.02163461538461538
$ws.Range("H19").Value = 0.202724358974359
$ws.Range("I19").Value = 0.07051282051282051
$ws.Range("J19").Value = 0.3501602564102564
$ws.Range("K19").Value = 0.1177884615384615
$ws.Range("M19").Value = 0.01682692307692308
$ws.Range("N19").Value = 0.001602564102564103
$ws.Range("O19").Value = 0.06570512820512821
$ws.Range("S19").Value = 0.1530448717948718
